# Update countries & provincias Spain
#
# Refresh the COVID-19 "Pais" sheet with a newer data pull (14 Apr 2020,
# 14:52 instead of 14:22) and fold in a few countries that were missing
# from the previous snapshot (Arabia Saudita, Guinea, Islas Caimanes,
# Burundi). Each affected row below is rewritten in full (country name +
# the 7 metric columns B:H) so the row ends up with the exact country
# label and figures from the refreshed dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 14 de Abril de 2020 a las 14:52"

# Row 33
$ws.Cells.Item(33,1).Value = "Dinamarca"
$ws.Cells.Item(33,2).Value = 6496
$ws.Cells.Item(33,3).Value = 178
$ws.Cells.Item(33,4).Value = 2515
$ws.Cells.Item(33,5).Value = 3682
$ws.Cells.Item(33,6).Value = 100
$ws.Cells.Item(33,7).Value = 14
$ws.Cells.Item(33,8).Value = 299

# Row 36
$ws.Cells.Item(36,1).Value = "Pakistan"
$ws.Cells.Item(36,2).Value = 5837
$ws.Cells.Item(36,3).Value = 341
$ws.Cells.Item(36,4).Value = 1378
$ws.Cells.Item(36,5).Value = 4363
$ws.Cells.Item(36,6).Value = 46
$ws.Cells.Item(36,7).Value = 3
$ws.Cells.Item(36,8).Value = 96

# Row 37
$ws.Cells.Item(37,1).Value = "Arabia Saudita"
$ws.Cells.Item(37,2).Value = 5369
$ws.Cells.Item(37,3).Value = 435
$ws.Cells.Item(37,4).Value = 889
$ws.Cells.Item(37,5).Value = 4407
$ws.Cells.Item(37,6).Value = 59
$ws.Cells.Item(37,7).Value = 8
$ws.Cells.Item(37,8).Value = 73

# Row 38
$ws.Cells.Item(38,1).Value = "Filipinas"
$ws.Cells.Item(38,2).Value = 5223
$ws.Cells.Item(38,3).Value = 291
$ws.Cells.Item(38,4).Value = 295
$ws.Cells.Item(38,5).Value = 4593
$ws.Cells.Item(38,6).Value = 1
$ws.Cells.Item(38,7).Value = 20
$ws.Cells.Item(38,8).Value = 335

# Row 39
$ws.Cells.Item(39,1).Value = "Mexico"
$ws.Cells.Item(39,2).Value = 5014
$ws.Cells.Item(39,3).Value = 353
$ws.Cells.Item(39,4).Value = 1964
$ws.Cells.Item(39,5).Value = 2718
$ws.Cells.Item(39,6).Value = 207
$ws.Cells.Item(39,7).Value = 36
$ws.Cells.Item(39,8).Value = 332

# Row 40
$ws.Cells.Item(40,1).Value = "Malasia"
$ws.Cells.Item(40,2).Value = 4987
$ws.Cells.Item(40,3).Value = 170
$ws.Cells.Item(40,4).Value = 2478
$ws.Cells.Item(40,5).Value = 2427
$ws.Cells.Item(40,6).Value = 60
$ws.Cells.Item(40,7).Value = 5
$ws.Cells.Item(40,8).Value = 82

# Row 52
$ws.Cells.Item(52,1).Value = "Colombia"
$ws.Cells.Item(52,2).Value = 2852
$ws.Cells.Item(52,3).Value = 0
$ws.Cells.Item(52,4).Value = 319
$ws.Cells.Item(52,5).Value = 2421
$ws.Cells.Item(52,6).Value = 106
$ws.Cells.Item(52,7).Value = 0
$ws.Cells.Item(52,8).Value = 112

# Row 70
$ws.Cells.Item(70,1).Value = "Kazajistan"
$ws.Cells.Item(70,2).Value = 1202
$ws.Cells.Item(70,3).Value = 111
$ws.Cells.Item(70,4).Value = 150
$ws.Cells.Item(70,5).Value = 1038
$ws.Cells.Item(70,6).Value = 21
$ws.Cells.Item(70,7).Value = 2
$ws.Cells.Item(70,8).Value = 14

# Row 105
$ws.Cells.Item(105,1).Value = "Guinea"
$ws.Cells.Item(105,2).Value = 363
$ws.Cells.Item(105,3).Value = 44
$ws.Cells.Item(105,4).Value = 31
$ws.Cells.Item(105,5).Value = 332
$ws.Cells.Item(105,6).Value = 0
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 0

# Row 106
$ws.Cells.Item(106,1).Value = "Bolivia"
$ws.Cells.Item(106,2).Value = 354
$ws.Cells.Item(106,3).Value = 24
$ws.Cells.Item(106,4).Value = 6
$ws.Cells.Item(106,5).Value = 320
$ws.Cells.Item(106,6).Value = 3
$ws.Cells.Item(106,7).Value = 1
$ws.Cells.Item(106,8).Value = 28

# Row 107
$ws.Cells.Item(107,1).Value = "Nigeria"
$ws.Cells.Item(107,2).Value = 343
$ws.Cells.Item(107,3).Value = 0
$ws.Cells.Item(107,4).Value = 91
$ws.Cells.Item(107,5).Value = 242
$ws.Cells.Item(107,6).Value = 2
$ws.Cells.Item(107,7).Value = 0
$ws.Cells.Item(107,8).Value = 10

# Row 108
$ws.Cells.Item(108,1).Value = "Mauricio"
$ws.Cells.Item(108,2).Value = 324
$ws.Cells.Item(108,3).Value = 0
$ws.Cells.Item(108,4).Value = 42
$ws.Cells.Item(108,5).Value = 273
$ws.Cells.Item(108,6).Value = 3
$ws.Cells.Item(108,7).Value = 0
$ws.Cells.Item(108,8).Value = 9

# Row 149
$ws.Cells.Item(149,1).Value = "Islas Caimanes"
$ws.Cells.Item(149,2).Value = 54
$ws.Cells.Item(149,3).Value = 0
$ws.Cells.Item(149,4).Value = 6
$ws.Cells.Item(149,5).Value = 47
$ws.Cells.Item(149,6).Value = 3
$ws.Cells.Item(149,7).Value = 0
$ws.Cells.Item(149,8).Value = 1

# Row 150
$ws.Cells.Item(150,1).Value = "Uganda"
$ws.Cells.Item(150,2).Value = 54
$ws.Cells.Item(150,3).Value = 0
$ws.Cells.Item(150,4).Value = 8
$ws.Cells.Item(150,5).Value = 46
$ws.Cells.Item(150,6).Value = 0
$ws.Cells.Item(150,7).Value = 0
$ws.Cells.Item(150,8).Value = 0

# Row 206
$ws.Cells.Item(206,1).Value = "Burundi"
$ws.Cells.Item(206,2).Value = 5
$ws.Cells.Item(206,3).Value = 0
$ws.Cells.Item(206,4).Value = 0
$ws.Cells.Item(206,5).Value = 4
$ws.Cells.Item(206,6).Value = 0
$ws.Cells.Item(206,7).Value = 0
$ws.Cells.Item(206,8).Value = 1

# Row 207
$ws.Cells.Item(207,1).Value = "Islas Malvinas"
$ws.Cells.Item(207,2).Value = 5
$ws.Cells.Item(207,3).Value = 0
$ws.Cells.Item(207,4).Value = 1
$ws.Cells.Item(207,5).Value = 4
$ws.Cells.Item(207,6).Value = 0
$ws.Cells.Item(207,7).Value = 0
$ws.Cells.Item(207,8).Value = 0
